{"js": "// Applies the \"two-digit division\" worksheet update:\n//   - the header date line: 2023-12-23 Saturday -> 2023-12-24 Sunday\n//   - each of the 25 division-problem cells gets a new problem/answer.\n// Every \"before\" string below is unique in the document, so an exact,\n// case-sensitive whole-document search/replace unambiguously targets the\n// correct run each time.\nconst replacements = [\n  [\"2023-12-23 Saturday\", \"2023-12-24 Sunday\"],\n  [\"98\u00f79=10, 8\", \"27\u00f79=3, 0\"],\n  [\"72\u00f77=10, 2\", \"14\u00f76=2, 2\"],\n  [\"59\u00f72=29, 1\", \"25\u00f73=8, 1\"],\n  [\"33\u00f74=8, 1\", \"49\u00f77=7, 0\"],\n  [\"82\u00f78=10, 2\", \"36\u00f73=12, 0\"],\n  [\"41\u00f75=8, 1\", \"84\u00f76=14, 0\"],\n  [\"90\u00f76=15, 0\", \"44\u00f77=6, 2\"],\n  [\"99\u00f75=19, 4\", \"89\u00f73=29, 2\"],\n  [\"74\u00f77=10, 4\", \"55\u00f76=9, 1\"],\n  [\"17\u00f72=8, 1\", \"24\u00f73=8, 0\"],\n  [\"78\u00f75=15, 3\", \"83\u00f74=20, 3\"],\n  [\"70\u00f78=8, 6\", \"83\u00f72=41, 1\"],\n  [\"44\u00f72=22, 0\", \"17\u00f78=2, 1\"],\n  [\"80\u00f76=13, 2\", \"77\u00f73=25, 2\"],\n  [\"18\u00f72=9, 0\", \"79\u00f78=9, 7\"],\n  [\"57\u00f74=14, 1\", \"29\u00f75=5, 4\"],\n  [\"54\u00f77=7, 5\", \"81\u00f75=16, 1\"],\n  [\"13\u00f79=1, 4\", \"50\u00f79=5, 5\"],\n  [\"96\u00f73=32, 0\", \"25\u00f78=3, 1\"],\n  [\"25\u00f72=12, 1\", \"26\u00f78=3, 2\"],\n  [\"59\u00f78=7, 3\", \"36\u00f78=4, 4\"],\n  [\"65\u00f77=9, 2\", \"57\u00f78=7, 1\"],\n  [\"68\u00f72=34, 0\", \"30\u00f75=6, 0\"],\n  [\"43\u00f76=7, 1\", \"55\u00f76=9, 1\"],\n  [\"95\u00f77=13, 4\", \"21\u00f75=4, 1\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + before);\n  }\n\n  for (const item of results.items) {\n    item.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Applies the \"two-digit division\" worksheet update:\n#   - the header date line: 2023-12-23 Saturday -> 2023-12-24 Sunday\n#   - each of the 25 division-problem cells gets a new problem/answer.\n# Every \"before\" string is unique in the document, so an exact,\n# case-sensitive whole-document Find/Replace unambiguously targets the\n# correct run each time.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2023-12-23 Saturday\", \"2023-12-24 Sunday\"),\n  @(\"98\u00f79=10, 8\", \"27\u00f79=3, 0\"),\n  @(\"72\u00f77=10, 2\", \"14\u00f76=2, 2\"),\n  @(\"59\u00f72=29, 1\", \"25\u00f73=8, 1\"),\n  @(\"33\u00f74=8, 1\", \"49\u00f77=7, 0\"),\n  @(\"82\u00f78=10, 2\", \"36\u00f73=12, 0\"),\n  @(\"41\u00f75=8, 1\", \"84\u00f76=14, 0\"),\n  @(\"90\u00f76=15, 0\", \"44\u00f77=6, 2\"),\n  @(\"99\u00f75=19, 4\", \"89\u00f73=29, 2\"),\n  @(\"74\u00f77=10, 4\", \"55\u00f76=9, 1\"),\n  @(\"17\u00f72=8, 1\", \"24\u00f73=8, 0\"),\n  @(\"78\u00f75=15, 3\", \"83\u00f74=20, 3\"),\n  @(\"70\u00f78=8, 6\", \"83\u00f72=41, 1\"),\n  @(\"44\u00f72=22, 0\", \"17\u00f78=2, 1\"),\n  @(\"80\u00f76=13, 2\", \"77\u00f73=25, 2\"),\n  @(\"18\u00f72=9, 0\", \"79\u00f78=9, 7\"),\n  @(\"57\u00f74=14, 1\", \"29\u00f75=5, 4\"),\n  @(\"54\u00f77=7, 5\", \"81\u00f75=16, 1\"),\n  @(\"13\u00f79=1, 4\", \"50\u00f79=5, 5\"),\n  @(\"96\u00f73=32, 0\", \"25\u00f78=3, 1\"),\n  @(\"25\u00f72=12, 1\", \"26\u00f78=3, 2\"),\n  @(\"59\u00f78=7, 3\", \"36\u00f78=4, 4\"),\n  @(\"65\u00f77=9, 2\", \"57\u00f78=7, 1\"),\n  @(\"68\u00f72=34, 0\", \"30\u00f75=6, 0\"),\n  @(\"43\u00f76=7, 1\", \"55\u00f76=9, 1\"),\n  @(\"95\u00f77=13, 4\", \"21\u00f75=4, 1\")\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n  $before = $pair[0]\n  $after = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $before\n  $find.Replacement.Text = $after\n  $find.Execute($find.Text, $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n}\n"}
